$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Harmonogram (schedule) fix: fill in the new K/L/M block for rows 71-75
# that was missing, mirroring the existing rows above (e.g. row 70),
# which record a date, a file name, and a line count.
$dateSerial = 45911

$rows = @(
    @{ Row = 71; Text = "harmonogram.component.ts";  Count = 5 },
    @{ Row = 72; Text = "login.component.css";       Count = 35 },
    @{ Row = 73; Text = "login.component.html";      Count = 4 },
    @{ Row = 74; Text = "register.component.css";    Count = 47 },
    @{ Row = 75; Text = "register.component.html";   Count = 4 }
)

foreach ($r in $rows) {
    # Copy K70 (which already carries the short-date number format) into
    # the target cell first so we inherit the existing style instead of
    # creating a brand-new number format, then overwrite the value.
    $ws.Range("K70").Copy($ws.Range("K$($r.Row)"))
    $ws.Range("K$($r.Row)").Value = $dateSerial

    $ws.Range("L$($r.Row)").Value = $r.Text
    $ws.Range("M$($r.Row)").Value = $r.Count
}

# Update the sheet's selected/active cell to match the recorded selection.
$ws.Range("M77").Select()
